# "code clean in student Basic info"
# - Update the test phone number and password text on the Login_credentials sheet
# - Drop the stale mailto: hyperlinks (and their Hyperlink cell styling) from A2/B2
# - Move the active selection off of B2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test login data
$ws.Range("A2").Value = 77122558691
$ws.Range("B2").Value = "Darshu@123"

# Remove the mailto hyperlinks that used to decorate A2/B2
$ws.Hyperlinks.Delete()

# Those cells were using the built-in "Hyperlink" cell style (underline, themed
# colour, border) -- clear that back to the default/no style now that the
# hyperlinks are gone, and drop the now-unused named style from the workbook.
$ws.Range("A2:B2").ClearFormats()
$wb.Styles("Hyperlink").Delete()

# Leave the selection somewhere other than the credential cell
$ws.Range("A6").Select()
